$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, pushing existing rows 88:185 down to 89:186.
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the latest weekly price record.
$ws.Range("A88").Value = 7
$ws.Range("B88").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C88").Value = "Ñuble"
$ws.Range("D88").Value = 44638
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112017
$ws.Range("G88").Value = "Apio"
$ws.Range("H88").Value = "Americana (o)"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 60
$ws.Range("K88").Value = 8500
$ws.Range("L88").Value = 9000
$ws.Range("M88").Value = 8750
$ws.Range("N88").Value = "$/docena de matas"
$ws.Range("O88").Value = "Provincia del Elquí"
$ws.Range("P88").Value = 1458
$ws.Range("Q88").Value = 6
$ws.Range("R88").Value = "Hortaliza"
